$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.810.63"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.088.57"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.11"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.81"
$ws.Range("E7").Value = "  +3.02%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.392"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0791"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.395.16"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.79"
$ws.Range("E13").Value = "  +3.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.25"
$ws.Range("E14").Value = "  +3.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.771"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.087.20"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.741.95"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.35"
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0837"
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.75"
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.27"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("E27").Value = "  +5.34%  "
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +2.39%  "
$ws.Range("E31").Value = "  +2.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.70"
$ws.Range("E32").Value = "  +3.00%  "
$ws.Range("E33").Value = "  +4.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0632"
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("E36").Value = "  +3.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.84"
$ws.Range("E37").Value = "  +3.17%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -3.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.100"
$ws.Range("E40").Value = "  +4.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.90"
$ws.Range("E41").Value = "  +2.10%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.38"
$ws.Range("E43").Value = "  +4.81%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0216"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.466.38"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("E47").Value = "  +4.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "16.05"
$ws.Range("E48").Value = "  +5.48%  "
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.279.99"
$ws.Range("E51").Value = "  +0.94%  "
